$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22. This pushes the existing row 22
# (date 44383 ...) down to row 23, and the existing row 23
# (date 44433 ...) down to row 24.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly entry.
$ws.Cells.Item(22, 1).Value = 11
$ws.Cells.Item(22, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(22, 3).Value = "Bíobío"
$ws.Cells.Item(22, 4).Value = 44509
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = 100114007
$ws.Cells.Item(22, 7).Value = "Jengibre"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 15000
$ws.Cells.Item(22, 12).Value = 16000
$ws.Cells.Item(22, 13).Value = 15500
$ws.Cells.Item(22, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(22, 15).Value = "Perú"
$ws.Cells.Item(22, 16).Value = 1192
$ws.Cells.Item(22, 17).Value = 13
$ws.Cells.Item(22, 18).Value = "Hortaliza"
